$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns A (date) and B (resource name) to fit the new report row
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws.Columns.Item(2).ColumnWidth = 52.666666666666664

# Append the new download-report row: date, resource name, count
$ws.Range("A3").Value = "1.8.2023 0:00:00"
$ws.Range("A3").HorizontalAlignment = -4108

$ws.Range("B3").Value = "C# 10 in a Nutshell The Definitive Reference Joseph Albahari"

$ws.Range("C3").Value = 1
